$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange

$part1 = '"There is a magical moment when a programmer presses the '
$part2 = 'run'
$part3 = ' button and the software begins to execute.  Somehow a program written in a '
$part4 = 'high-level language is '
$part5 = 'running on a computer that is capable only of shuffling bits.  Here we reveal the wizardry that makes that moment possible." ' + [char]0x2013 + ' Jeremy '
$part6 = 'Siek'

$quote = $part1 + $part2 + $part3 + $part4 + $part5 + $part6

$tr.Text = $quote
[void]$tr.InsertAfter("`r")
$tr.Font.Size = 24

$pos1 = 1
$pos2 = $pos1 + $part1.Length
$pos3 = $pos2 + $part2.Length
$pos4 = $pos3 + $part3.Length
$pos5 = $pos4 + $part4.Length
$pos6 = $pos5 + $part5.Length

$tr.Characters($pos1, $part1.Length).Font.Size = 24
$tr.Characters($pos2, $part2.Length).Font.Italic = $true
$tr.Characters($pos3, $part3.Length).Font.Size = 24
$tr.Characters($pos4, $part4.Length).Font.Size = 24
$tr.Characters($pos5, $part5.Length).Font.Size = 24
$tr.Characters($pos6, $part6.Length).Font.Size = 24
